$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update province ranking: Asturias and Murcia moved up due to new data ---

# Row 11: was Alacant/Alicante 603/12/541/50 -> now Asturias 662/30/610/22
$ws.Range("A11").Value = "Asturias"
$ws.Range("B11").Value = 662
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 610
$ws.Range("E11").Value = 22

# Row 12: was Asturias 545/12/467/10 -> now Alacant/Alicante 603/12/541/50
$ws.Range("A12").Value = "Alacant/Alicante"
$ws.Range("B12").Value = 603
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 541
$ws.Range("E12").Value = 50

# Row 20: was Gipuzkoa/Guipuzcoa 380/283/365/15 -> now Murcia 385/1/381/3
$ws.Range("A20").Value = "Murcia"
$ws.Range("B20").Value = 385
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 381
$ws.Range("E20").Value = 3

# Row 21: was Granada 374/0/357/17 -> now Gipuzkoa/Guipuzcoa 380/283/365/15
$ws.Range("A21").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B21").Value = 380
$ws.Range("C21").Value = 283
$ws.Range("D21").Value = 365
$ws.Range("E21").Value = 15

# Row 22: was Sevilla 351/1/345/5 -> now Granada 374/0/357/17
$ws.Range("A22").Value = "Granada"
$ws.Range("B22").Value = 374
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 357
$ws.Range("E22").Value = 17

# Row 23: was Valladolid 349/17/318/14 -> now Sevilla 351/1/345/5
$ws.Range("A23").Value = "Sevilla"
$ws.Range("B23").Value = 351
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 345
$ws.Range("E23").Value = 5

# Row 24: was Cantabria 347/11/330/6 -> now Valladolid 349/17/318/14
$ws.Range("A24").Value = "Valladolid"
$ws.Range("B24").Value = 349
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 318
$ws.Range("E24").Value = 14

# Row 25: was Murcia 345/1/213/3 -> now Cantabria 347/11/330/6
$ws.Range("A25").Value = "Cantabria"
$ws.Range("B25").Value = 347
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 330
$ws.Range("E25").Value = 6

# --- Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 21:16"
